$d = $word.ActiveDocument

# Replace the entire body content with the target OOXML fragment (the
# trailing <w:sectPr> is left untouched since we only touch $d.Content,
# which spans up to -- but not including -- the section properties).
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:noSpellErr="1">
  <w:pPr><w:jc w:val="center"/></w:pPr>
  <w:r>
    <w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:eastAsia="Calibri" w:cs="Calibri"/><w:b w:val="1"/><w:bCs w:val="1"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr>
    <w:t xml:space="preserve">SI304/ST366 – Engenharia de Software II                                                                                     </w:t>
  </w:r>
  <w:r>
    <w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:eastAsia="Calibri" w:cs="Calibri"/><w:b w:val="1"/><w:bCs w:val="1"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr>
    <w:t>1 º Sem/2016</w:t>
  </w:r>
</w:p>
<w:p w14:noSpellErr="1"/>
<w:p w14:noSpellErr="1">
  <w:r>
    <w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:eastAsia="Calibri" w:cs="Calibri"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr>
    <w:t>Caroline Resende Silveira – ra: 165921</w:t>
  </w:r>
</w:p>
<w:p w14:noSpellErr="1">
  <w:r>
    <w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:eastAsia="Calibri" w:cs="Calibri"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr>
    <w:t>Karina Sayuri Hagiwara – ra: 171419</w:t>
  </w:r>
</w:p>
<w:p w14:noSpellErr="1"/>
<w:p w14:noSpellErr="1">
  <w:r>
    <w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:eastAsia="Calibri" w:cs="Calibri"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr>
    <w:t>1. Elicitação de requisitos</w:t>
  </w:r>
</w:p>
<w:p w14:noSpellErr="1">
  <w:r><w:br/></w:r>
</w:p>
<w:p w14:noSpellErr="1">
  <w:r>
    <w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:eastAsia="Calibri" w:cs="Calibri"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr>
    <w:t>2. Requisitos conscientes</w:t>
  </w:r>
</w:p>
<w:p w14:noSpellErr="1">
  <w:r>
    <w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:eastAsia="Calibri" w:cs="Calibri"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr>
    <w:t xml:space="preserve">  o O sistema deverá verificar se um número é impar ou par quando um número inteiro positivo é fornecido.</w:t>
  </w:r>
</w:p>
<w:p w14:noSpellErr="1">
  <w:r>
    <w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:eastAsia="Calibri" w:cs="Calibri"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr>
    <w:t xml:space="preserve">  o O sistema deverá informar ao usuário se o número é negativo ou positivo quando um número inteiro é fornecido pelo usuário.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:r>
    <w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:eastAsia="Calibri" w:cs="Calibri"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr>
    <w:t xml:space="preserve">  o Assim que iniciado, o sistema deverá mostrar ao usuário um menu e verificar que operação ele deseja realizar – verificar se um </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:eastAsia="Calibri" w:cs="Calibri"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr>
    <w:t>número é</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:eastAsia="Calibri" w:cs="Calibri"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr>
    <w:t xml:space="preserve"> par ou impar ou verificar se um numero é positivo ou negativo ou se ele deseja sair do programa.</w:t>
  </w:r>
</w:p>
<w:p w14:noSpellErr="1">
  <w:r>
    <w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:eastAsia="Calibri" w:cs="Calibri"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr>
    <w:t xml:space="preserve">  o O sistema deverá verificar qual a opção que o usuário deseja realizar.</w:t>
  </w:r>
</w:p>
<w:p w14:noSpellErr="1">
  <w:r><w:br/></w:r>
</w:p>
<w:p w:rsidR="0313C27A" w:rsidP="0313C27A" w:rsidRDefault="0313C27A" w14:noSpellErr="1" w14:paraId="36CC4C49" w14:textId="6471385E">
  <w:pPr><w:pStyle w:val="Normal"/></w:pPr>
</w:p>
'@

$d.Content.InsertXML($xml) | Out-Null

Write-Host "Applied target XML fragment."
